# RF001 - Autenticar Usuario : update Step/Expected-Result text for the
# "wrong password" transition across TC1..TC6 (swap between "Usuario do
# Sistema preenche os campos..." and "...seleciona um nome de usuario
# sugerido..." as step 2, with the matching Expected Result).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$preenche   = "Usuario do Sistema preenche os campos e clica no botao entrar"
$seleciona  = "Usuario do Sistema seleciona um nome de usuario sugerido, digita a senha e clica no botao entrar"
$alertaCAS  = "SYSTEM alerta que o CAS (sistema de autorizacao login-senha) esta fora do ar"
$alertaUsr  = "SYSTEM alerta que o nome de usuario e/ou senha estao incorretos"
$alertaTJ   = "SYSTEM alerta que o TJSeg (sistema que fornece as permissoes de acesso e escrita) esta fora do ar"

# TC1 (rows 10-13): step 2
$ws.Range("B11").Value = $seleciona
$ws.Range("D11").Value = $alertaUsr

# TC2 (rows 20-23): step 2
$ws.Range("B21").Value = $preenche
$ws.Range("D21").Value = $alertaCAS

# TC3 (rows 30-33): step 2
$ws.Range("B31").Value = $preenche
$ws.Range("D31").Value = $alertaUsr

# TC4 (rows 40-43): step 2 and step 3
$ws.Range("B41").Value = $preenche
$ws.Range("D41").Value = $alertaTJ
$ws.Range("B42").Value = $seleciona

# TC5 (rows 50-53): step 2
$ws.Range("B51").Value = $seleciona
$ws.Range("D51").Value = $alertaCAS

# TC6 (rows 60-63): step 2 and step 3
$ws.Range("B61").Value = $seleciona
$ws.Range("B62").Value = $preenche
